$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure numeric-looking price strings in column D stay as text (matches source data formatting)
$rngD = $ws.Range("D2:D51")
$rngD.NumberFormat = "@"

# Apply updated cryptocurrency price and volume figures
$ws.Range('D2').Value = '29.479.85'
$ws.Range('E2').Value = '  +0.78%  '
$ws.Range('D3').Value = '1.919.92'
$ws.Range('E3').Value = '  +1.57%  '
$ws.Range('E4').Value = '  +0.83%  '
$ws.Range('D5').Value = '326.11'
$ws.Range('E5').Value = '  +1.05%  '
$ws.Range('D7').Value = '0.4831'
$ws.Range('E7').Value = '  +2.55%  '
$ws.Range('D8').Value = '0.4087'
$ws.Range('E8').Value = '  +1.27%  '
$ws.Range('D9').Value = '0.08244'
$ws.Range('E9').Value = '  +2.89%  '
$ws.Range('D10').Value = '1.024'
$ws.Range('E10').Value = '  +3.00%  '
$ws.Range('D11').Value = '23.51'
$ws.Range('E11').Value = '  +3.19%  '
$ws.Range('D12').Value = '1.939.68'
$ws.Range('E12').Value = '  +1.28%  '
$ws.Range('D13').Value = '6.049'
$ws.Range('E13').Value = '  +2.07%  '
$ws.Range('D14').Value = '7.243'
$ws.Range('E14').Value = '  +3.08%  '
$ws.Range('D15').Value = '91.32'
$ws.Range('E15').Value = '  +2.19%  '
$ws.Range('D16').Value = '0.06800'
$ws.Range('E16').Value = '  +2.61%  '
$ws.Range('D17').Value = '1.008'
$ws.Range('E17').Value = '  +0.72%  '
$ws.Range('D18').Value = '0.00001041'
$ws.Range('E18').Value = '  +1.87%  '
$ws.Range('D19').Value = '17.78'
$ws.Range('E19').Value = '  +2.00%  '
$ws.Range('E20').Value = '  +0.64%  '
$ws.Range('D21').Value = '29.518.33'
$ws.Range('E21').Value = '  +0.93%  '
$ws.Range('D22').Value = '5.645'
$ws.Range('E22').Value = '  +2.85%  '
$ws.Range('D23').Value = '11.79'
$ws.Range('E23').Value = '  +1.04%  '
$ws.Range('E24').Value = '  +0.93%  '
$ws.Range('D25').Value = '2.168.47'
$ws.Range('E25').Value = '  +4.18%  '
$ws.Range('D26').Value = '6.629'
$ws.Range('E26').Value = '  +10.62%  '
$ws.Range('D27').Value = '156.75'
$ws.Range('E27').Value = '  +1.06%  '
$ws.Range('D28').Value = '20.12'
$ws.Range('E28').Value = '  +2.51%  '
$ws.Range('D29').Value = '2.120'
$ws.Range('E29').Value = '  +1.79%  '
$ws.Range('D30').Value = '120.52'
$ws.Range('E30').Value = '  +3.07%  '
$ws.Range('D31').Value = '1.025'
$ws.Range('E31').Value = '  +0.35%  '
$ws.Range('D32').Value = '0.09579'
$ws.Range('E32').Value = '  +1.83%  '
$ws.Range('D33').Value = '5.531'
$ws.Range('E33').Value = '  +3.40%  '
$ws.Range('D34').Value = '3.563'
$ws.Range('E34').Value = '  +0.71%  '
$ws.Range('E35').Value = '  +0.26%  '
$ws.Range('E36').Value = '  +1.89%  '
$ws.Range('D37').Value = '0.06136'
$ws.Range('E37').Value = '  +1.50%  '
$ws.Range('D38').Value = '1.178'
$ws.Range('E38').Value = '  +0.63%  '
$ws.Range('D39').Value = '0.5990'
$ws.Range('E39').Value = '  +2.99%  '
$ws.Range('D40').Value = '8.049'
$ws.Range('E40').Value = '  +0.79%  '
$ws.Range('E41').Value = '  +7.62%  '
$ws.Range('E42').Value = '  +1.42%  '
$ws.Range('D43').Value = '1.280'
$ws.Range('E43').Value = '  +0.69%  '
$ws.Range('D44').Value = '2.402'
$ws.Range('E44').Value = '  +0.73%  '
$ws.Range('D45').Value = '0.07607'
$ws.Range('E45').Value = '  -1.32%  '
$ws.Range('D46').Value = '12.43'
$ws.Range('E46').Value = '  +2.35%  '
$ws.Range('D47').Value = '0.5583'
$ws.Range('E47').Value = '  +2.00%  '
$ws.Range('D48').Value = '1.960'
$ws.Range('E48').Value = '  +3.01%  '
$ws.Range('D49').Value = '117.74'
$ws.Range('E49').Value = '  +3.95%  '
$ws.Range('E50').Value = '  +4.63%  '
$ws.Range('D51').Value = '72.77'
$ws.Range('E51').Value = '  +2.31%  '

# Restore default (unformatted) style so saved cells match original styling
$rngD.Style = "Normal"
